$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.355.74'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.58%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.876.79'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.55%  '

# Row 4
$ws.Range('E4').Value = '  +0.07%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.09'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.07%  '

# Row 6
$ws.Range('E6').Value = '  +0.10%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4771'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.73%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2881'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.96%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06514'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.66%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.24'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.77%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07763'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.09%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.880.12'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.59%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '96.22'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.31%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.7325'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.96%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.116'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.25%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '275.28'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.40%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.339.36'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.56%  '

# Row 18
$ws.Range('E18').Value = '  -2.20%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007531'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.46%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.001'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.12%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.123.44'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.08%  '

# Row 22
$ws.Range('E22').Value = '  +0.09%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.221'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.63%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.156'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.08%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.212'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.73%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.06'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.44%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.92'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.24%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.951'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.77%  '

# Row 29
$ws.Range('E29').Value = '  +0.29%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09952'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.24%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.505'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.14%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.303'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.32%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.076'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.80%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04738'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.22%  '

# Row 35
$ws.Range('E35').Value = '  -0.79%  '

# Row 36
$ws.Range('E36').Value = '  -1.02%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.719'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.03%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01850'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.35%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.745'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.93%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.278'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.27%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8417'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.61%  '

# Row 42
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.906'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.36%  '

# Row 43
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.000'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.09%  '

# Row 44
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4159'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.42%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '69.12'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.49%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '101.73'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.12%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.267'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.03%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.078'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.44%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '35.10'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.45%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '909.94'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -6.24%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05586'
$ws.Range('D51').Style = 'Normal'
